$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '23.233.25'
$ws.Range("E2").Value = '  +12.99%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.628.56'
$ws.Range("E3").Value = '  +10.60%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9906'
$ws.Range("E4").Value = '  -1.80%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '304.51'
$ws.Range("E5").Value = '  +9.87%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.9833'
$ws.Range("E6").Value = '  +2.55%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3691'
$ws.Range("E7").Value = '  +3.52%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3430'
$ws.Range("E8").Value = '  +11.72%  '
$ws.Range("B9").Value = 'Polygon'
$ws.Range("C9").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '1.160'
$ws.Range("E9").Value = '  +6.36%  '
$ws.Range("B10").Value = 'OKB'
$ws.Range("C10").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '42.39'
$ws.Range("E10").Value = '  +7.57%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07121'
$ws.Range("E11").Value = '  +7.20%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.9863'
$ws.Range("E12").Value = '  -1.64%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '20.32'
$ws.Range("E13").Value = '  +12.05%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.928'
$ws.Range("E14").Value = '  +8.61%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.657'
$ws.Range("E15").Value = '  +7.76%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.00001085'
$ws.Range("E16").Value = '  +6.15%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.622.89'
$ws.Range("E17").Value = '  +10.27%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.9821'
$ws.Range("E18").Value = '  +2.47%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06681'
$ws.Range("E19").Value = '  +12.11%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '78.66'
$ws.Range("E20").Value = '  +13.72%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '16.22'
$ws.Range("E21").Value = '  +11.54%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.027'
$ws.Range("E22").Value = '  +9.78%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '11.77'
$ws.Range("E23").Value = '  +4.24%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '23.188.36'
$ws.Range("E24").Value = '  +12.74%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.373'
$ws.Range("E25").Value = '  +4.33%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '3.426'
$ws.Range("E26").Value = '  -7.23%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.633'
$ws.Range("E27").Value = '  +26.20%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '150.28'
$ws.Range("E28").Value = '  +3.67%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '19.39'
$ws.Range("E29").Value = '  +13.04%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.799.29'
$ws.Range("E30").Value = '  +10.28%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '124.85'
$ws.Range("E31").Value = '  +9.30%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.064'
$ws.Range("E32").Value = '  +5.24%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '6.115'
$ws.Range("E33").Value = '  +23.86%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.9849'
$ws.Range("E34").Value = '  +24.12%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.08290'
$ws.Range("E35").Value = '  +4.50%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.673'
$ws.Range("E36").Value = '  +15.66%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '8.802'
$ws.Range("E37").Value = '  +21.03%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '12.00'
$ws.Range("E38").Value = '  +16.34%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '5.227'
$ws.Range("E39").Value = '  +11.04%  '
$ws.Range("B40").Value = 'Hedera'
$ws.Range("C40").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.06228'
$ws.Range("E40").Value = '  +8.52%  '
$ws.Range("B41").Value = 'TrustWalletToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.263'
$ws.Range("E41").Value = '  +2.02%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.02249'
$ws.Range("E42").Value = '  +10.62%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.2029'
$ws.Range("E43").Value = '  +9.26%  '
$ws.Range("B44").Value = 'TheSandbox'
$ws.Range("C44").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.5972'
$ws.Range("E44").Value = '  +13.80%  '
$ws.Range("B45").Value = 'Frax'
$ws.Range("C45").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.9825'
$ws.Range("E45").Value = '  +2.45%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.811'
$ws.Range("E46").Value = '  +8.65%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '13.09'
$ws.Range("E47").Value = '  +8.83%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.5811'
$ws.Range("E48").Value = '  +12.13%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '126.91'
$ws.Range("E49").Value = '  +7.02%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.994'
$ws.Range("E50").Value = '  +10.63%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.06955'
$ws.Range("E51").Value = '  +8.05%  '
